$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that *looks* numeric (e.g. "176.28") while keeping it
# stored as literal text, matching the source data's inlineStr cells. Plain
# `.Value = "176.28"` would be auto-coerced to a number by Excel's normal
# cell-entry parsing, so instead we stage the text in a scratch cell that is
# explicitly formatted as Text ("@"), copy it, and paste-special *values only*
# into the destination. PasteSpecial(xlPasteValues) carries over the text
# value without carrying over the scratch cell's Text number format, so the
# destination cell keeps its original (General) style.
$stage = $ws.Range("ZZ1")
$stage.NumberFormat = "@"

function Set-TextValue($addr, $text) {
    $stage.Value = $text
    $stage.Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}

# --- Cells whose new text is not number-like: safe to assign directly ---
$ws.Range("D2").Value = "66.858.19"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "3.085.11"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("E5").Value = "  -0.87%  "
$ws.Range("E6").Value = "  +4.61%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.083.61"
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("E10").Value = "  +0.47%  "
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("E12").Value = "  -1.17%  "
$ws.Range("E13").Value = "  -0.52%  "
$ws.Range("E14").Value = "  -0.68%  "
$ws.Range("E15").Value = "  +0.87%  "
$ws.Range("D16").Value = "3.599.98"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").Value = "66.851.98"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("E18").Value = "  -0.40%  "
$ws.Range("D19").Value = "3.084.50"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("E20").Value = "  -1.14%  "
$ws.Range("E21").Value = "  -1.65%  "
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("E23").Value = "  -0.53%  "
$ws.Range("E24").Value = "  +0.54%  "
$ws.Range("E25").Value = "  -1.34%  "
$ws.Range("E26").Value = "  +0.90%  "
$ws.Range("E27").Value = "  +0.83%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("E31").Value = "  -1.49%  "
$ws.Range("E32").Value = "  +0.28%  "
$ws.Range("E33").Value = "  -0.24%  "
$ws.Range("D34").Value = "0.0₃0929"
$ws.Range("E34").Value = "  +2.74%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("E36").Value = "  -1.37%  "
$ws.Range("E37").Value = "  -2.52%  "
$ws.Range("E38").Value = "  -0.48%  "
$ws.Range("E39").Value = "  +3.03%  "
$ws.Range("E40").Value = "  -0.55%  "
$ws.Range("E41").Value = "  +0.46%  "
$ws.Range("E42").Value = "  -0.33%  "
$ws.Range("E43").Value = "  -1.16%  "
$ws.Range("E44").Value = "  +9.38%  "
$ws.Range("D45").Value = "2.798.96"
$ws.Range("E45").Value = "  +1.01%  "
$ws.Range("E46").Value = "  -0.59%  "
$ws.Range("E47").Value = "  -0.67%  "
$ws.Range("E48").Value = "  -0.59%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("E50").Value = "  +4.98%  "
$ws.Range("E51").Value = "  +6.96%  "

# --- Cells whose new text parses as a plain number: must be forced to text ---
Set-TextValue "D6" "176.28"
Set-TextValue "D9" "0.514"
Set-TextValue "D10" "6.42"
Set-TextValue "D14" "35.83"
Set-TextValue "D18" "6.98"
Set-TextValue "D20" "16.47"
Set-TextValue "D21" "483.68"
Set-TextValue "D22" "7.69"
Set-TextValue "D25" "12.69"
Set-TextValue "D27" "10.30"
Set-TextValue "D29" "7.84"
Set-TextValue "D30" "2.29"
Set-TextValue "D38" "46.65"
Set-TextValue "D39" "0.312"
Set-TextValue "D44" "2.67"
Set-TextValue "D46" "369.67"
Set-TextValue "D48" "134.45"
Set-TextValue "D50" "25.76"
Set-TextValue "D51" "2.31"

$stage.Clear()
